$d = $word.ActiveDocument

# 1. Title text change
$d.Content.Find.Execute("Serial Transmit of Temperature", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Measurement and Transmission of Propeller Speed", 2)

# 2. Course code - merge three runs into a single run's text (gramStart/gramEnd proofErr runs get removed)
$d.Content.Find.Execute("EN.605.715.81.FA19 - Software Development for Real-Time Systems", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EN.605.715.81.FA19 - Software Development for Real-Time Systems", 2)

# 3. RPM over time. - merge runs
$d.Content.Find.Execute("RPM over time.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RPM over time.", 2)
